# BOT; UPDATE DATA
# Adds the 2020-05-23 (Excel serial 43974) daily row to the three data
# sheets ("all", "kobe", "other"), pushing the trailing footnote/label
# row down by one, and updates the active selection on each sheet.

$wb = $excel.ActiveWorkbook

# ---- Sheet "all" : new row 46 (date 43974) -----------------------------
$ws1 = $wb.Worksheets.Item("all")
$ws1.Activate()
$ws1.Rows.Item(46).Insert()
$ws1.Range("A46").Value = 43974
$ws1.Range("B46").Value = 285
$ws1.Range("C46").Value = 282
$ws1.Range("D46").Value = 36
$ws1.Range("E46").Value = 32
$ws1.Range("F46").Value = 4
$ws1.Range("G46").Value = 12
$ws1.Range("H46").Value = 234
$ws1.Range("C49").Select() | Out-Null

# ---- Sheet "kobe" : new row 101 (date 43974) ----------------------------
$ws2 = $wb.Worksheets.Item("kobe")
$ws2.Activate()
$ws2.Rows.Item(101).Insert()
$ws2.Range("A101").Value = 43974
$ws2.Range("B101").Value = 0
$ws2.Range("C101").Value = 3010
$ws2.Range("D101").Value = 0
$ws2.Range("E101").Value = 285
$ws2.Range("F101").Value = 31
$ws2.Range("G101").Value = 28
$ws2.Range("H101").Value = 3
$ws2.Range("I101").Value = 12
$ws2.Range("J101").Value = 225
$ws2.Range("K101").Select() | Out-Null

# ---- Sheet "other" : new row 76 (date 43974) ----------------------------
$ws3 = $wb.Worksheets.Item("other")
$ws3.Activate()
$ws3.Rows.Item(76).Insert()
$ws3.Range("A76").Value = 43974
$ws3.Range("B76").Value = 0
$ws3.Range("C76").Value = 14
$ws3.Range("D76").Value = 5
$ws3.Range("E76").Value = 4
$ws3.Range("F76").Value = 1
$ws3.Range("G76").Value = 0
$ws3.Range("H76").Value = 9
$ws3.Range("F77").Select() | Out-Null

# Restore the originally active sheet/tab ("all").
$ws1.Activate()
